$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.164.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.57%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.587.21'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.32%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.20'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.15%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.20'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.72%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.20%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.53%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.598.60'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.66%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.14%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.82%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.72%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.21%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.048.86'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.19%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.938.84'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.25%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.58'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.17%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.577.83'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.13%  '

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.32%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '338.44'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.52%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.82%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.07'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.47%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.23%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.14'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.29%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.13%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.01%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.18%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.02'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.28%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.03%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.88%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.72%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.48%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.71'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.61%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.45'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.28%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.82%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.90%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.80'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.94%  '

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.45%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.825'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.41%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.816'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.27%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.35%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.33%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '272.65'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.86%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.78'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.15%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.592'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.90%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0952'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.42%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0516'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.47%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.42'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.83%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.967.82'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.76%  '

$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0220'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.86%  '

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.51'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.82%  '
